# Trade #2 closed at 2026-02-17 13:33:15 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet: update capital / P&L / trade counters ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.97   # Current Capital
$summary.Range("B4").Value = -0.03     # Total P&L $
$summary.Range("B5").Value = -0.3      # Total P&L %
$summary.Range("B6").Value = 2         # Total Trades
$summary.Range("B8").Value = 2         # Losing Trades

# --- Strategy Status sheet: update MarketMaking strategy row ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.97      # Capital
$status.Range("D4").Value = 2          # Trades
$status.Range("E4").Value = -0.03      # P&L $
$status.Range("F4").Value = -0.03      # P&L %

# --- Add the new trade row (#2) to both "All Trades" and "MarketMaking" sheets ---
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(3, 1).Value = 2
    # Force the date column to be stored as literal text (matching column B's
    # existing "Date" column which holds text like "2026-02-17", not a real
    # date) - otherwise Excel's COM Value setter auto-parses it into a date
    # serial number.
    $ws.Cells.Item(3, 2).NumberFormat = "@"
    $ws.Cells.Item(3, 2).Value = "2026-02-17"
    $ws.Cells.Item(3, 3).Value = "13:33:09"
    $ws.Cells.Item(3, 4).Value = "MarketMaking"
    $ws.Cells.Item(3, 5).Value = "UP"
    $ws.Cells.Item(3, 6).Value = 0.16
    $ws.Cells.Item(3, 7).Value = 0.14538
    $ws.Cells.Item(3, 8).Value = "CLOSED"
    $ws.Cells.Item(3, 9).Value = -9.137600000000001
    $ws.Cells.Item(3, 10).Value = -0.01
    $ws.Cells.Item(3, 11).Value = 99.97
    $ws.Cells.Item(3, 12).Value = 0
    $ws.Cells.Item(3, 13).Value = 0
    $ws.Cells.Item(3, 14).Value = 0.6
    $ws.Cells.Item(3, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(3, 16).Value = "early_exit"
    $ws.Cells.Item(3, 17).Value = 0.12
}
